$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") rows 2-70, per regenerated save_data (uses K instead of Strike#)
$gValues = @{
    2 = 1
    3 = 2
    4 = 3
    5 = 3
    6 = 1
    7 = 2
    8 = 1
    9 = 2
    10 = 1
    11 = 3
    12 = 2
    13 = 3
    14 = 0
    15 = 1
    16 = 1
    17 = 1
    18 = 1
    19 = 2
    20 = 2
    21 = 0
    22 = 2
    23 = 2
    24 = 1
    25 = 2
    26 = 1
    27 = 1
    28 = 1
    29 = 1
    30 = 0
    31 = 3
    32 = 1
    33 = 2
    34 = 3
    35 = 2
    36 = 1
    37 = 2
    38 = 1
    39 = 3
    40 = 2
    41 = 2
    42 = 0
    43 = 1
    44 = 0
    45 = 1
    46 = 2
    47 = 3
    48 = 1
    49 = 3
    50 = 0
    51 = 1
    52 = 1
    53 = 2
    54 = 1
    55 = 2
    56 = 3
    57 = 4
    58 = 1
    59 = 2
    60 = 2
    61 = 3
    62 = 3
    63 = 3
    64 = 0
    65 = 1
    66 = 2
    67 = 1
    68 = 4
    69 = 1
    70 = 2
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}

